$d = $word.ActiveDocument

# The target paragraph is the empty paragraph that sits between the
# first and second tables in the document (just after the bookmarkEnd
# that closes table 1). Locate it robustly by matching its start
# position against the end of table 1, instead of relying on a
# hard-coded paragraph index.
$table1End = $d.Tables(1).Range.End

$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $candidate = $d.Paragraphs($i)
    if ($candidate.Range.Start -eq $table1End) {
        $target = $candidate
        break
    }
}

# Select the empty paragraph (its range is just the paragraph mark)
# and apply single underline formatting through the Selection object,
# which correctly stores the formatting on the paragraph mark's run
# properties (w:pPr/w:rPr/w:u) even though the paragraph has no text.
$target.Range.Select()
$word.Selection.Underline = 1
